$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 139. Existing rows 139..158 shift down to 140..159,
# preserving their current values (including styles).
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new weekly record.
$ws.Cells.Item(139, 1).Value = 8
$ws.Cells.Item(139, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(139, 3).Value = "Coquimbo"
$ws.Cells.Item(139, 4).Value = 45131
$ws.Cells.Item(139, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(139, 5).Value = 4
$ws.Cells.Item(139, 6).Value = 100114007
$ws.Cells.Item(139, 7).Value = "Jengibre"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 360
$ws.Cells.Item(139, 11).Value = 17000
$ws.Cells.Item(139, 12).Value = 18000
$ws.Cells.Item(139, 13).Value = 17500
$ws.Cells.Item(139, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(139, 15).Value = "Perú"
$ws.Cells.Item(139, 16).Value = 1346
$ws.Cells.Item(139, 17).Value = 13
$ws.Cells.Item(139, 18).Value = "Hortaliza"
